# Auto-generated script to update leve profit calculation values
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets (scheduled price refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 308.17648
$ws.Range("I18").Value = 308.17648
$ws.Range("K18").Value = 308.17648
$ws.Range("M18").Value = -24.17648000000003
$ws.Range("H43").Value = 3556.7222
$ws.Range("J43").Value = 3554.4285
$ws.Range("L43").Value = 3554.4285
$ws.Range("N43").Value = -3692.4285
$ws.Range("H62").Value = 3339.2
$ws.Range("J62").Value = 2209.1428
$ws.Range("L62").Value = 2209.1428
$ws.Range("N62").Value = -3457.1428
$ws.Range("H65").Value = 3339.2
$ws.Range("J65").Value = 2209.1428
$ws.Range("L65").Value = 11045.714
$ws.Range("N65").Value = -17285.714
$ws.Range("H116").Value = 4313.375
$ws.Range("J116").Value = 4882.533
$ws.Range("L116").Value = 4882.533
$ws.Range("N116").Value = -11766.533
$ws.Range("H137").Value = 14827.903
$ws.Range("J137").Value = 21052.945
$ws.Range("L137").Value = 63158.835
$ws.Range("N137").Value = -68258.83499999999
$ws.Range("H138").Value = 1739.95
$ws.Range("J138").Value = 2497
$ws.Range("L138").Value = 7491
$ws.Range("N138").Value = -17771

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6203.5
$ws.Range("I45").Value = 7968.6
$ws.Range("K45").Value = 7968.6
$ws.Range("M45").Value = -7591.6
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 16777.445
$ws.Range("I20").Value = 23053.334
$ws.Range("J20").Value = 8932.583000000001
$ws.Range("K20").Value = 23053.334
$ws.Range("L20").Value = 8932.583000000001
$ws.Range("M20").Value = -22806.334
$ws.Range("N20").Value = -9426.583000000001
$ws.Range("H64").Value = 236.9
$ws.Range("I64").Value = 119.333336
$ws.Range("J64").Value = 287.2857
$ws.Range("K64").Value = 119.333336
$ws.Range("L64").Value = 287.2857
$ws.Range("M64").Value = 105.666664
$ws.Range("N64").Value = -737.2857
$ws.Range("H67").Value = 236.9
$ws.Range("I67").Value = 119.333336
$ws.Range("J67").Value = 287.2857
$ws.Range("K67").Value = 119.333336
$ws.Range("L67").Value = 287.2857
$ws.Range("M67").Value = 660.666664
$ws.Range("N67").Value = -1847.2857
$ws.Range("H134").Value = 8196.615
$ws.Range("I134").Value = 4235.697
$ws.Range("K134").Value = 12707.091
$ws.Range("M134").Value = -10172.091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 30601
$ws.Range("J92").Value = 30601
$ws.Range("L92").Value = 30601
$ws.Range("N92").Value = -35593
$ws.Range("H93").Value = 8999.5
$ws.Range("I93").Value = 5333
$ws.Range("K93").Value = 5333
$ws.Range("M93").Value = -3461
$ws.Range("H95").Value = 6930.4287
$ws.Range("J95").Value = 6930.4287
$ws.Range("L95").Value = 6930.4287
$ws.Range("N95").Value = -12422.4287
$ws.Range("H96").Value = 14561.6
$ws.Range("J96").Value = 14561.6
$ws.Range("L96").Value = 14561.6
$ws.Range("N96").Value = -20053.6
$ws.Range("H99").Value = 9063.929
$ws.Range("I99").Value = 4374.375
$ws.Range("K99").Value = 4374.375
$ws.Range("M99").Value = -2876.375
$ws.Range("H126").Value = 9063.929
$ws.Range("I126").Value = 4374.375
$ws.Range("K126").Value = 13123.125
$ws.Range("M126").Value = -10653.125
$ws.Range("H140").Value = 100000
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360
$ws.Range("H141").Value = 427307.06
$ws.Range("I141").Value = 324999
$ws.Range("J141").Value = 457999.5
$ws.Range("K141").Value = 324999
$ws.Range("L141").Value = 457999.5
$ws.Range("M141").Value = -319819
$ws.Range("N141").Value = -468359.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 28312482
$ws.Range("I4").Value = 30993008
$ws.Range("K4").Value = 92979024
$ws.Range("M4").Value = -92978912
$ws.Range("H12").Value = 608.6429000000001
$ws.Range("I12").Value = 26.25
$ws.Range("J12").Value = 841.6
$ws.Range("K12").Value = 78.75
$ws.Range("L12").Value = 2524.8
$ws.Range("M12").Value = 94.25
$ws.Range("N12").Value = -2870.8
$ws.Range("H102").Value = 8489.429
$ws.Range("I102").Value = 3142
$ws.Range("K102").Value = 9426
$ws.Range("M102").Value = -6992
$ws.Range("H106").Value = 12500
$ws.Range("J106").Value = 12500
$ws.Range("L106").Value = 37500
$ws.Range("N106").Value = -39392
$ws.Range("H116").Value = 1539.2
$ws.Range("I116").Value = 1539.2
$ws.Range("K116").Value = 4617.6
$ws.Range("M116").Value = -1175.6
$ws.Range("H131").Value = 2235.4
$ws.Range("J131").Value = 2274.638
$ws.Range("L131").Value = 6823.914
$ws.Range("N131").Value = -16903.914

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H70").Value = 11317.833
$ws.Range("I70").Value = 8852.25
$ws.Range("K70").Value = 8852.25
$ws.Range("M70").Value = -8582.25
$ws.Range("H73").Value = 11317.833
$ws.Range("I73").Value = 8852.25
$ws.Range("K73").Value = 8852.25
$ws.Range("M73").Value = -7916.25
$ws.Range("H80").Value = 4670.1
$ws.Range("J80").Value = 5023.75
$ws.Range("L80").Value = 5023.75
$ws.Range("N80").Value = -7019.75
$ws.Range("H83").Value = 4670.1
$ws.Range("J83").Value = 5023.75
$ws.Range("L83").Value = 25118.75
$ws.Range("N83").Value = -35102.75
$ws.Range("H92").Value = 28898.75
$ws.Range("I92").Value = 9200
$ws.Range("J92").Value = 31712.857
$ws.Range("K92").Value = 9200
$ws.Range("L92").Value = 31712.857
$ws.Range("M92").Value = -7328
$ws.Range("N92").Value = -35456.857
$ws.Range("H97").Value = 740
$ws.Range("I97").Value = 528
$ws.Range("J97").Value = 1058
$ws.Range("K97").Value = 528
$ws.Range("L97").Value = 1058
$ws.Range("M97").Value = -32
$ws.Range("N97").Value = -2050
$ws.Range("H122").Value = 8111.2417
$ws.Range("I122").Value = 6223.1353
$ws.Range("K122").Value = 18669.4059
$ws.Range("M122").Value = -16219.4059

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4302
$ws.Range("I40").Value = 3232.8572
$ws.Range("K40").Value = 3232.8572
$ws.Range("M40").Value = -3096.8572
$ws.Range("H42").Value = 35465.668
$ws.Range("I42").Value = 39399
$ws.Range("J42").Value = 33499
$ws.Range("K42").Value = 39399
$ws.Range("L42").Value = 33499
$ws.Range("M42").Value = -38836
$ws.Range("N42").Value = -34625
$ws.Range("H44").Value = 19800
$ws.Range("J44").Value = 19800
$ws.Range("L44").Value = 19800
$ws.Range("N44").Value = -20712
$ws.Range("H49").Value = 35465.668
$ws.Range("I49").Value = 39399
$ws.Range("J49").Value = 33499
$ws.Range("K49").Value = 39399
$ws.Range("L49").Value = 33499
$ws.Range("M49").Value = -39252
$ws.Range("N49").Value = -33793
$ws.Range("H68").Value = 3275.3333
$ws.Range("I68").Value = 2742.5
$ws.Range("K68").Value = 2742.5
$ws.Range("M68").Value = -1993.5
$ws.Range("H71").Value = 3275.3333
$ws.Range("I71").Value = 2742.5
$ws.Range("K71").Value = 13712.5
$ws.Range("M71").Value = -9968.5
$ws.Range("H94").Value = 15000
$ws.Range("J94").Value = 15000
$ws.Range("L94").Value = 15000
$ws.Range("N94").Value = -16352
$ws.Range("H100").Value = 4134.6
$ws.Range("I100").Value = 4769.4
$ws.Range("J100").Value = 3499.8
$ws.Range("K100").Value = 4769.4
$ws.Range("L100").Value = 3499.8
$ws.Range("M100").Value = -4228.4
$ws.Range("N100").Value = -4581.8
$ws.Range("H132").Value = 5731.28
$ws.Range("I132").Value = 3775.9375
$ws.Range("K132").Value = 11327.8125
$ws.Range("M132").Value = -8797.8125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 54351.332
$ws.Range("I49").Value = 31527.5
$ws.Range("J49").Value = 99999
$ws.Range("K49").Value = 31527.5
$ws.Range("L49").Value = 99999
$ws.Range("M49").Value = -31297.5
$ws.Range("N49").Value = -100459
$ws.Range("H122").Value = 5571.25
$ws.Range("I122").Value = 4234.9585
$ws.Range("J122").Value = 8243.833000000001
$ws.Range("K122").Value = 12704.8755
$ws.Range("L122").Value = 24731.499
$ws.Range("M122").Value = -10254.8755
$ws.Range("N122").Value = -29631.499
